$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data range (previously A1:E1) before laying out the new table
$ws.Range("A1:E1").ClearContents()

# Header row
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "pierdoli"
$ws.Range("D1").Value = "smiedzi"

# Data rows - format column A as text first so the numeric-looking IDs
# ("1", "2") are stored as text, not numbers, then restore Normal style
$ws.Range("A2:A3").NumberFormat = "@"

$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "debil@gmail.com"
$ws.Range("C2").Value = "no debil no"
$ws.Range("D2").Value = "N/A"

$ws.Range("A3").Value = "2"
$ws.Range("B3").Value = "idiota@gmail.com"
$ws.Range("C3").Value = "N/A"
$ws.Range("D3").Value = "no idiota no"

$ws.Range("A2:A3").Style = "Normal"
